$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.710906360849698
$ws.Range("E2").Value = 0.710906360849698

# Row 3
$ws.Range("D3").Value = 0.001364785198323819
$ws.Range("E3").Value = 0.001364785198323819

# Row 4
$ws.Range("D4").Value = [double]"3.234733729112554E-05"
$ws.Range("E4").Value = [double]"3.234733729112554E-05"

# Row 5
$ws.Range("D5").Value = 0.1050741683980147
$ws.Range("E5").Value = 0.1050741683980147

# Row 6
$ws.Range("D6").Value = 0.961422305147399
$ws.Range("E6").Value = 0.961422305147399

# Row 7
$ws.Range("D7").Value = 0.9999999743748171
$ws.Range("E7").Value = [double]"2.562518286453042E-08"

# Row 8
$ws.Range("D8").Value = 0.5265769948507761
$ws.Range("E8").Value = 0.4734230051492239

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.4808591486035605
$ws.Range("E9").Value = 0.5191408513964395

# Row 10
$ws.Range("D10").Value = 0.8921866206081822
$ws.Range("E10").Value = 0.1078133793918178

# Row 11
$ws.Range("D11").Value = 0.5069167267431223
$ws.Range("E11").Value = 0.4930832732568777
$ws.Range("F11").Value = 0.6775525212287903
$ws.Range("G11").Value = 0.7

# Row 12
$ws.Range("D12").Value = 0.6300775253661615
$ws.Range("E12").Value = 0.6300775253661615

# Row 13
$ws.Range("D13").Value = 0.5521401895959576
$ws.Range("E13").Value = 0.5521401895959576

# Row 14
$ws.Range("D14").Value = [double]"3.675055887364767E-06"
$ws.Range("E14").Value = [double]"3.675055887364767E-06"

# Row 15
$ws.Range("D15").Value = 0.03464964306534804
$ws.Range("E15").Value = 0.03464964306534804

# Row 16
$ws.Range("D16").Value = 0.9990986540763646
$ws.Range("E16").Value = 0.9990986540763646

# Row 17
$ws.Range("D17").Value = 0.999999997471656
$ws.Range("E17").Value = [double]"2.528343978625003E-09"

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.0342812851180604
$ws.Range("E18").Value = 0.9657187148819396

# Row 19
$ws.Range("D19").Value = 0.7490431635262348
$ws.Range("E19").Value = 0.2509568364737652

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"8.735485527413459E-05"
$ws.Range("E20").Value = 0.9999126451447259

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.1238257725780879
$ws.Range("E21").Value = 0.8761742274219122
$ws.Range("F21").Value = 2.394115209579468
$ws.Range("G21").Value = 0.4
